$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4101058542728424
$ws.Range("B1").Value = 0.5166788697242737
$ws.Range("C1").Value = 0.7507784962654114
$ws.Range("D1").Value = 3.398271799087524
$ws.Range("E1").Value = 5.901803493499756
